$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set "none" into column G (Type of Specialty) for every data row (2-41)
# that doesn't already have a value there.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq $null) {
        $cell.Value = "none"
    }
}

# Update frozen pane top-left cell and active selection
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Range("A11").Select() | Out-Null
$ws.Range("G25:G41").Select() | Out-Null
